$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1: "Ordenar Datos" title shape -- the two runs ("     " and
# "Ordenar Datos") that share identical formatting get merged into a single
# run with the combined text "     Ordenar Datos".
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(1)
$tr1 = $shp1.TextFrame.TextRange

$secondRun = $tr1.Characters(6, $tr1.Length - 5)
$secondRun.Delete()

$firstRun = $tr1.Characters(1, 5)
$firstRun.InsertAfter("Ordenar Datos")

# Re-normalize the fill color (no visual change, white stays white).
$tr1b = $shp1.TextFrame.TextRange
$tr1b.Font.Color.RGB = 16777215

# ---------------------------------------------------------------------------
# Slide 4: practice-instructions shape -- drop the final bullet paragraph
# about "JERARQUIA.EQV" (that exercise step was removed).
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(1)
$tr4 = $shp4.TextFrame.TextRange

$paras4 = $tr4.Paragraphs()
$lastParaIdx = $paras4.Count
$lastPara = $tr4.Paragraphs($lastParaIdx, 1)
$lastPara.Delete()

# The delete above can leave a residual empty trailing paragraph behind;
# remove it too so the shape ends up with exactly six paragraphs.
$residual = $tr4.Paragraphs($lastParaIdx, 1)
if ($residual.Length -eq 0) {
    $residual.Delete()
}
